# Planet sheet refresh: shift the 8-day window forward (new dates + new
# measurements) and move the active selection from H13 to H14, matching a
# newer export of the same tracking template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 1-12 (row 13 is unchanged in the source diff).
$data = @(
    @(44175, 44176, 44177, 44178, 44179, 44180, 44181, 44182),
    @(259, 260, 261, 262, 263, 264, 265, 266),
    @(206, 220, 235, 250, 265, 279, 293, 307),
    @(254, 256, 257, 259, 260, 262, 263, 265),
    @(234, 235, 236, 238, 239, 240, 241, 243),
    @(20, 20, 20, 20, 21, 21, 21, 22),
    @(336, 336, 337, 337, 337, 337, 338, 338),
    @(298, 298, 299, 299, 299, 299, 299, 300),
    @(299, 299, 300, 300, 300, 300, 300, 300),
    @(37, 37, 37, 37, 37, 37, 37, 37),
    @(348, 348, 348, 348, 348, 348, 348, 348),
    @(294, 294, 294, 294, 294, 294, 294, 294)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

# Move the selected cell from H13 to H14.
$ws.Range("H14").Select()
